$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.365.15"
$ws.Range("E2").Value = "  +6.03%  "
$ws.Range("D3").Value = "2.362.57"
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "110.14"
$ws.Range("E5").Value = "  +3.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "309.52"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.616"
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.29"
$ws.Range("E10").Value = "  +2.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0918"
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.51"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.986"
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").Value = "2.726.83"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "2.368.14"
$ws.Range("E17").Value = "  +2.52%  "
$ws.Range("D18").Value = "45.363.33"
$ws.Range("E18").Value = "  +6.01%  "
$ws.Range("E19").Value = "  -2.49%  "
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.16"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.55"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "259.82"
$ws.Range("E24").Value = "  -2.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.28"
$ws.Range("E25").Value = "  +1.49%  "
$ws.Range("E26").Value = "  -0.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.11"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.39"
$ws.Range("E28").Value = "  -5.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.35"
$ws.Range("E29").Value = "  +2.32%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0966"
$ws.Range("E30").Value = "  +10.77%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.41"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "38.01"
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "170.52"
$ws.Range("E33").Value = "  +2.79%  "
$ws.Range("E34").Value = "  +7.29%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.82"
$ws.Range("E36").Value = "  +3.80%  "
$ws.Range("E37").Value = "  +3.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.98"
$ws.Range("E38").Value = "  +6.80%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.91"
$ws.Range("E40").Value = "  +7.43%  "
$ws.Range("E41").Value = "  +8.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.53"
$ws.Range("E42").Value = "  -4.17%  "
$ws.Range("E43").Value = "  +0.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.99"
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.99"
$ws.Range("E45").Value = "  +5.53%  "
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "82.02"
$ws.Range("E47").Value = "  +7.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.39"
$ws.Range("E48").Value = "  +5.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "112.97"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.55"
$ws.Range("E50").Value = "  +6.15%  "
$ws.Range("D51").Value = "1.628.33"
$ws.Range("E51").Value = "  -3.97%  "
